$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 741 (shifts the existing row 741..782 down to 742..783,
# preserving every cell's value/type, exactly like the diff's record insertion).
$ws.Rows.Item(741).Insert()

# Populate the newly inserted row with the new record.
# The date/weekday columns in this sheet are stored as literal text (not real
# Excel dates), so a leading apostrophe keeps "2026/02/02" from being
# auto-converted into a date serial number; resetting the style back to
# "Normal" afterwards drops the quotePrefix formatting COM tacks on, so the
# cell ends up with the same (default) style as its neighbours.
$ws.Range("A741").Value = "'2026/02/02"
$ws.Range("A741").Style = "Normal"

$ws.Range("B741").Value = "'月"
$ws.Range("B741").Style = "Normal"

$ws.Range("C741").Value = 8
$ws.Range("D741").Value = 22
